$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores dotted/subscript-formatted numeric strings
# (not real numbers), so force text formatting before writing the
# refreshed values to avoid Excel re-interpreting them as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.557.20"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "1.842.90"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "231.95"
$ws.Range("E5").Value = "  +3.56%  "
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  +3.67%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "44.01"
$ws.Range("E8").Value = "  +12.39%  "
$ws.Range("E9").Value = "  +8.40%  "
$ws.Range("D10").Value = "0.0699"
$ws.Range("E10").Value = "  +4.98%  "
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("D12").Value = "2.108.28"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("D13").Value = "1.856.87"
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").Value = "11.27"
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").Value = "0.673"
$ws.Range("E15").Value = "  +7.63%  "
$ws.Range("D16").Value = "4.72"
$ws.Range("E16").Value = "  +8.65%  "
$ws.Range("D17").Value = "35.516.50"
$ws.Range("E17").Value = "  +3.25%  "
$ws.Range("D18").Value = "70.47"
$ws.Range("E18").Value = "  +3.87%  "
$ws.Range("D19").Value = "0.0₃0801"
$ws.Range("E19").Value = "  +5.42%  "
$ws.Range("D20").Value = "244.50"
$ws.Range("E20").Value = "  +2.51%  "
$ws.Range("E21").Value = "  +9.21%  "
$ws.Range("D22").Value = "4.65"
$ws.Range("E22").Value = "  +14.38%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  +3.94%  "
$ws.Range("D25").Value = "171.22"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").Value = "8.01"
$ws.Range("E26").Value = "  +5.05%  "
$ws.Range("D27").Value = "17.83"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("E29").Value = "  +28.15%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "3.300.67"
$ws.Range("E31").Value = "  +35.85%  "
$ws.Range("D32").Value = "0.0552"
$ws.Range("E32").Value = "  +7.88%  "
$ws.Range("E33").Value = "  +7.06%  "
$ws.Range("D34").Value = "3.94"
$ws.Range("E34").Value = "  +5.83%  "
$ws.Range("D35").Value = "1.84"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").Value = "95.05"
$ws.Range("E36").Value = "  +16.89%  "
$ws.Range("E37").Value = "  +8.70%  "
$ws.Range("D38").Value = "1.14"
$ws.Range("E38").Value = "  +7.88%  "
$ws.Range("D39").Value = "1.348.56"
$ws.Range("E39").Value = "  +3.93%  "
$ws.Range("D40").Value = "0.0196"
$ws.Range("E40").Value = "  +6.26%  "
$ws.Range("D41").Value = "15.38"
$ws.Range("E41").Value = "  +10.01%  "
$ws.Range("E42").Value = "  +6.68%  "
$ws.Range("E43").Value = "  +8.23%  "
$ws.Range("E44").Value = "  +3.55%  "
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("E47").Value = "  +9.90%  "
$ws.Range("D48").Value = "0.0519"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "2.013.40"
$ws.Range("E49").Value = "  +2.79%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "102.48"
$ws.Range("E51").Value = "  +1.16%  "
